# Add 5 new rows (363-367) to the "liste référence" sheet, describing
# hydrogen-production related variables/colors, matching the structure of
# the existing rows (columns: variable / label_fr / label_en / color).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data -----------------------------------------------------------
# English variable / label_en text (column A and C share the same text)
$englishLabels = @(
    "Biomass gasification (with CCS)",
    "Biomass gasification (without CCS)",
    "Electrolyzers",
    "Autothermal reforming (with CCS)",
    "Steam methane reforming (without CCS)"
)

# French label (column B)
$frenchLabels = @(
    "Gazéification de biomasse (avec CSC)",
    "Gazéification de biomasse (sans CSC)",
    "Électrolyseur",
    "Reformage à la vapeur authothermique (avec CSC)",
    "Reformage à la vapeur (sans CSC)"
)

# Color hex code shown in column D, used to derive the fill (as BGR decimal
# values expected by the Interior.Color COM property) for each row.
$colorHex = @("#e8d4b7", "#60613b", "#77dba4", "#d29d31", "#f6b4a4")
$colorBgr = @(12047592, 3891552, 10804087, 3251666, 10794230)

$startRow = 363

# --- Column A: write first so the shared-string table picks up the five
# English strings before anything else (matches indices 826-830). -------
# Row 365 keeps the default look-and-feel (not reset to "Normal") on
# columns A/C, matching the source workbook; the others are reset so they
# carry no explicit style (as in the source file).
for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $cell = $ws.Range("A$row")
    $cell.Value = $englishLabels[$i]
    if ($row -ne 365) {
        $cell.Style = "Normal"
    }
}

# --- Column B: French labels (indices 831-835). -------------------------
for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $ws.Range("B$row").Value = $frenchLabels[$i]
}

# --- Column C: reuses the same shared strings as column A. --------------
for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $cell = $ws.Range("C$row")
    $cell.Value = $englishLabels[$i]
    if ($row -ne 365) {
        $cell.Style = "Normal"
    }
}

# --- Column D: color codes, with a fill color matching the hex value. ---
# The last new row (367) uses an explicit black font color (set before the
# fill so the resulting format matches the style already used elsewhere in
# the sheet for that particular fill, rather than creating a brand-new,
# unused style record).
for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $cell = $ws.Range("D$row")
    $cell.Value = $colorHex[$i]
    if ($row -eq 367) {
        $cell.Font.Color = 0
    }
    $cell.Interior.Color = $colorBgr[$i]
}

# --- Update the view to point at the newly added rows, mirroring what a
# user would see after scrolling down to/selecting the last entered cell.
$ws.Range("D367").Select()

